# ---------------------------------------------------------------------------
# Edit script: applies the changes described by the XML diff to before.docx
#
#   1. Removes the stray <w:bookmarkStart.../><w:bookmarkEnd/> pair
#      (id="0" name="_GoBack") that sits in the first, empty paragraph.
#   2. Renumbers the five legacy VML picture shape ids down by one
#      (_x0000_i1026..1030 -> _x0000_i1025..1029) and truncates their
#      o:title attributes from "...at 00.58" to "...at 00".
#   3. Appends a new empty paragraph followed by a paragraph containing
#      "Selaamla be yeni update" (dotless-i), with a fresh _GoBack bookmark
#      placed at the very end of the document (its usual "last edit" spot).
#
# The body is rebuilt from a known-good literal copy of the original
# word/document.xml content (rather than re-derived via the WordOpenXML
# property, which normalizes/merges runs and injects w14:paraId -- losing
# fidelity with the source) and written back with Range.InsertXML, which
# replaces the exact range's contents in place.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- Step 1: drop the old _GoBack bookmark; it is re-created at its new
#     end-of-document location in step 3 below. ---------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- The original body markup (verbatim), used as the editing template. ---
$body = @'
<w:p w:rsidR="009C29AE" w:rsidRDefault="00EA264E"><w:r><w:t xml:space="preserve">                         </w:t></w:r><w:r w:rsidR="009C29AE"><w:t xml:space="preserve">             </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p w:rsidR="009C29AE" w:rsidRDefault="009C29AE"><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C853BE"><w:pict><v:shapetype id="_x0000_t75" coordsize="21600,21600" o:spt="75" o:preferrelative="t" path="m@4@5l@4@11@9@11@9@5xe" filled="f" stroked="f"><v:stroke joinstyle="miter"/><v:formulas><v:f eqn="if lineDrawn pixelLineWidth 0"/><v:f eqn="sum @0 1 0"/><v:f eqn="sum 0 0 @1"/><v:f eqn="prod @2 1 2"/><v:f eqn="prod @3 21600 pixelWidth"/><v:f eqn="prod @3 21600 pixelHeight"/><v:f eqn="sum @0 0 1"/><v:f eqn="prod @6 1 2"/><v:f eqn="prod @7 21600 pixelWidth"/><v:f eqn="sum @8 21600 0"/><v:f eqn="prod @7 21600 pixelHeight"/><v:f eqn="sum @10 21600 0"/></v:formulas><v:path o:extrusionok="f" gradientshapeok="t" o:connecttype="rect"/><o:lock v:ext="edit" aspectratio="t"/></v:shapetype><v:shape id="_x0000_i1026" type="#_x0000_t75" style="width:151.8pt;height:189.6pt"><v:imagedata r:id="rId4" o:title="WhatsApp Image 2025-06-25 at 00.58"/></v:shape></w:pict></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00C853BE"><w:pict><v:shape id="_x0000_i1027" type="#_x0000_t75" style="width:142.8pt;height:197.4pt"><v:imagedata r:id="rId5" o:title="WhatsApp Image 2025-06-25 at 00.58"/></v:shape></w:pict></w:r><w:r w:rsidR="00C853BE"><w:pict><v:shape id="_x0000_i1028" type="#_x0000_t75" style="width:156.6pt;height:190.2pt"><v:imagedata r:id="rId6" o:title="WhatsApp Image 2025-06-25 at 00.58"/></v:shape></w:pict></w:r></w:p><w:p w:rsidR="00D95DDE" w:rsidRDefault="00C853BE"><w:r><w:pict><v:shape id="_x0000_i1029" type="#_x0000_t75" style="width:137.4pt;height:173.4pt"><v:imagedata r:id="rId7" o:title="WhatsApp Image 2025-06-25 at 00.58"/></v:shape></w:pict></w:r><w:r w:rsidR="00EA264E"><w:t>t</w:t></w:r><w:r><w:pict><v:shape id="_x0000_i1030" type="#_x0000_t75" style="width:165pt;height:205.8pt"><v:imagedata r:id="rId8" o:title="WhatsApp Image 2025-06-25 at 00.58"/></v:shape></w:pict></w:r><w:r w:rsidR="009C29AE"><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0"><wp:extent cx="1898073" cy="2403673"/><wp:effectExtent l="0" t="0" r="6985" b="0"/><wp:docPr id="2" name="Resim 2" descr="C:\Users\ramazan\AppData\Local\Microsoft\Windows\INetCache\Content.Word\WhatsApp Image 2025-06-25 at 00.58.31 (1).jpeg"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 42" descr="C:\Users\ramazan\AppData\Local\Microsoft\Windows\INetCache\Content.Word\WhatsApp Image 2025-06-25 at 00.58.31 (1).jpeg"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId9" cstate="print"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="1954350" cy="2474941"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>
'@

# --- Step 2a: remove the leftover _GoBack bookmark markup from the text.
$body = $body -replace '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>', ''

# --- Step 2b: renumber the VML shape ids (_x0000_iNNNN -> _x0000_i(NNNN-1)).
#     The scriptblock/delegate form of -replace isn't available in this
#     restricted runtime, so the five known ids are substituted explicitly
#     via unique placeholder tokens first, so chained replacements never
#     collide with each other (e.g. 1027->1026 must not then be caught by
#     a later 1026->1025 rule).
$body = $body -replace '_x0000_i1026', '@@SHAPE1025@@'
$body = $body -replace '_x0000_i1027', '@@SHAPE1026@@'
$body = $body -replace '_x0000_i1028', '@@SHAPE1027@@'
$body = $body -replace '_x0000_i1029', '@@SHAPE1028@@'
$body = $body -replace '_x0000_i1030', '@@SHAPE1029@@'
$body = $body -replace '@@SHAPE1025@@', '_x0000_i1025'
$body = $body -replace '@@SHAPE1026@@', '_x0000_i1026'
$body = $body -replace '@@SHAPE1027@@', '_x0000_i1027'
$body = $body -replace '@@SHAPE1028@@', '_x0000_i1028'
$body = $body -replace '@@SHAPE1029@@', '_x0000_i1029'

# --- Step 2c: truncate the picture titles.
$body = $body -replace 'WhatsApp Image 2025-06-25 at 00\.58(?=")', 'WhatsApp Image 2025-06-25 at 00'

# --- Step 3: append the new trailing paragraphs (empty paragraph, then the
#     "Selaamla be yeni update" paragraph carrying the relocated _GoBack
#     bookmark). ----------------------------------------------------------
$dotlessI = [string][char]0x0131
$newParagraphs = '<w:p/>' +
    '<w:p>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>Selaamla</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> be </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>yen' + $dotlessI + '</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> update</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>'

$body = $body + $newParagraphs

# --- Declare the namespaces this fragment relies on (w, r, v, o, wp, a,
#     pic, a14) on its first element so InsertXML can parse it standalone.
$nsDecls = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"' +
    ' xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"' +
    ' xmlns:v="urn:schemas-microsoft-com:vml"' +
    ' xmlns:o="urn:schemas-microsoft-com:office:office"' +
    ' xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing"' +
    ' xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"' +
    ' xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"' +
    ' xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"'

$firstElemEnd = $body.IndexOf(">")
$body = $body.Substring(0, $firstElemEnd) + $nsDecls + $body.Substring($firstElemEnd)

# --- Replace the whole document content with the rebuilt body. ------------
$full = $d.Content
$full.InsertXML($body)
